$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - update values
$ws.Range("B3").Value = 0.01757020049857353
$ws.Range("C3").Value = 0.01911226359807061
$ws.Range("D3").Value = 61282315316520.35

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.02004270100233599
$ws.Range("C4").Value = 0.02202844305256954
$ws.Range("D4").Value = 0.082088466149235

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 21879283314028.84
$ws.Range("C5").Value = 4396964490890.669
$ws.Range("D5").Value = 164074473178293.1
